$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 8827
$ws.Range("J51").Value = 7320.3335
$ws.Range("L51").Value = 7320.3335
$ws.Range("N51").Value = -8288.333500000001
$ws.Range("H106").Value = 2398.7932
$ws.Range("I106").Value = 2317.8845
$ws.Range("K106").Value = 2317.8845
$ws.Range("M106").Value = -1686.8845
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = ""
$ws.Range("H113").Value = 1350
$ws.Range("J113").Value = 1500
$ws.Range("L113").Value = 1500
$ws.Range("N113").Value = -8008
$ws.Range("H135").Value = 793.4
$ws.Range("I135").Value = 666
$ws.Range("J135").Value = 878.3333
$ws.Range("K135").Value = 5994
$ws.Range("L135").Value = 7904.9997
$ws.Range("M135").Value = -3459
$ws.Range("N135").Value = -12974.9997
$ws.Range("H137").Value = 7237.4736
$ws.Range("I137").Value = 2212.3333
$ws.Range("J137").Value = 11760.1
$ws.Range("K137").Value = 6636.999899999999
$ws.Range("L137").Value = 35280.3
$ws.Range("M137").Value = -4086.999899999999
$ws.Range("N137").Value = -40380.3

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 48612
$ws.Range("J7").Value = 48612
$ws.Range("L7").Value = 48612
$ws.Range("N7").Value = -48840
$ws.Range("H26").Value = 6004
$ws.Range("I26").Value = 4502
$ws.Range("J26").Value = 7506
$ws.Range("K26").Value = 4502
$ws.Range("L26").Value = 7506
$ws.Range("M26").Value = -4172
$ws.Range("N26").Value = -8166
$ws.Range("H32").Value = 3821.4238
$ws.Range("I32").Value = 3187.0574
$ws.Range("J32").Value = 14859.4
$ws.Range("K32").Value = 3187.0574
$ws.Range("L32").Value = 14859.4
$ws.Range("M32").Value = -2900.0574
$ws.Range("N32").Value = -15433.4
$ws.Range("H45").Value = 55668.375
$ws.Range("I45").Value = 420420
$ws.Range("K45").Value = 420420
$ws.Range("M45").Value = -420043
$ws.Range("H61").Value = 8097.533
$ws.Range("I61").Value = 2507.3333
$ws.Range("K61").Value = 2507.3333
$ws.Range("M61").Value = -2295.3333
$ws.Range("H74").Value = 187564.36
$ws.Range("I74").Value = 279510.4
$ws.Range("J74").Value = 3672.3
$ws.Range("K74").Value = 279510.4
$ws.Range("L74").Value = 3672.3
$ws.Range("M74").Value = -278636.4
$ws.Range("N74").Value = -5420.3
$ws.Range("H77").Value = 187564.36
$ws.Range("I77").Value = 279510.4
$ws.Range("J77").Value = 3672.3
$ws.Range("K77").Value = 1397552
$ws.Range("L77").Value = 18361.5
$ws.Range("M77").Value = -1393184
$ws.Range("N77").Value = -27097.5
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").Value = ""
$ws.Range("H136").Value = 8097.533
$ws.Range("I136").Value = 2507.3333
$ws.Range("K136").Value = 7521.999899999999
$ws.Range("M136").Value = -4971.999899999999
$ws.Range("H141").Value = 137497.5
$ws.Range("J141").Value = 137497.5
$ws.Range("L141").Value = 137497.5
$ws.Range("N141").Value = -147857.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 25005568
$ws.Range("I20").Value = 31255752
$ws.Range("K20").Value = 31255752
$ws.Range("M20").Value = -31255505
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").Value = ""
$ws.Range("H107").Value = 1544.3889
$ws.Range("I107").Value = 1406.25
$ws.Range("J107").Value = 2649.5
$ws.Range("K107").Value = 1406.25
$ws.Range("L107").Value = 2649.5
$ws.Range("M107").Value = 513.75
$ws.Range("N107").Value = -6489.5
$ws.Range("H132").Value = 83000
$ws.Range("J132").Value = 83000
$ws.Range("L132").Value = 83000
$ws.Range("N132").Value = -93120

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 9591.666999999999
$ws.Range("I32").Value = 9554.637000000001
$ws.Range("J32").Value = 9999
$ws.Range("K32").Value = 9554.637000000001
$ws.Range("L32").Value = 9999
$ws.Range("M32").Value = -9238.637000000001
$ws.Range("N32").Value = -10631
$ws.Range("H86").Value = 4366.625
$ws.Range("I86").Value = 4156
$ws.Range("K86").Value = 4156
$ws.Range("M86").Value = -3033
$ws.Range("H89").Value = 4366.625
$ws.Range("I89").Value = 4156
$ws.Range("K89").Value = 20780
$ws.Range("M89").Value = -15164
$ws.Range("H94").Value = 2458.4167
$ws.Range("I94").Value = 1878.6
$ws.Range("J94").Value = 2872.5715
$ws.Range("K94").Value = 1878.6
$ws.Range("L94").Value = 2872.5715
$ws.Range("M94").Value = -1427.6
$ws.Range("N94").Value = -3774.5715
$ws.Range("H106").Value = 363499.5
$ws.Range("J106").Value = 363499.5
$ws.Range("L106").Value = 363499.5
$ws.Range("N106").Value = -366023.5
$ws.Range("H132").Value = 4942.0835
$ws.Range("J132").Value = 5388.636
$ws.Range("L132").Value = 16165.908
$ws.Range("N132").Value = -21225.908
$ws.Range("H134").Value = 2665.6086
$ws.Range("I134").Value = 2376.5293
$ws.Range("J134").Value = 3484.6667
$ws.Range("K134").Value = 7129.5879
$ws.Range("L134").Value = 10454.0001
$ws.Range("M134").Value = -4594.5879
$ws.Range("N134").Value = -15524.0001
$ws.Range("H141").Value = 360237.5
$ws.Range("J141").Value = 360237.5
$ws.Range("L141").Value = 360237.5
$ws.Range("N141").Value = -370597.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 2527
$ws.Range("J108").Value = 5499.6665
$ws.Range("L108").Value = 16498.9995
$ws.Range("N108").Value = -22258.9995

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 125002490
$ws.Range("I80").Value = 166668530
$ws.Range("K80").Value = 166668530
$ws.Range("M80").Value = -166667532
$ws.Range("H83").Value = 125002490
$ws.Range("I83").Value = 166668530
$ws.Range("K83").Value = 833342650
$ws.Range("M83").Value = -833337658
$ws.Range("H102").Value = 1709.75
$ws.Range("I102").Value = 1496.25
$ws.Range("K102").Value = 1496.25
$ws.Range("M102").Value = 125.75
$ws.Range("H136").Value = 47333.332
$ws.Range("J136").Value = 47333.332
$ws.Range("L136").Value = 141999.996
$ws.Range("N136").Value = -147099.996

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").Value = ""
$ws.Range("H61").Value = 1843.9565
$ws.Range("I61").Value = 1932
$ws.Range("K61").Value = 1932
$ws.Range("M61").Value = -1730
$ws.Range("H113").Value = 1843.9565
$ws.Range("I113").Value = 1932
$ws.Range("K113").Value = 1932
$ws.Range("M113").Value = 238
$ws.Range("H122").Value = 5071.9
$ws.Range("I122").Value = 5411.5835
$ws.Range("J122").Value = 4562.375
$ws.Range("K122").Value = 16234.7505
$ws.Range("L122").Value = 13687.125
$ws.Range("M122").Value = -13784.7505
$ws.Range("N122").Value = -18587.125
$ws.Range("H136").Value = 9083.083000000001
$ws.Range("I136").Value = 8499.666999999999
$ws.Range("K136").Value = 25499.001
$ws.Range("M136").Value = -22949.001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 43279
$ws.Range("J75").Value = 48749
$ws.Range("L75").Value = 48749
$ws.Range("N75").Value = -50621
$ws.Range("H78").Value = 43279
$ws.Range("J78").Value = 48749
$ws.Range("L78").Value = 146247
$ws.Range("N78").Value = -155607
$ws.Range("H88").Value = 183
$ws.Range("I88").Value = 183
$ws.Range("K88").Value = 183
$ws.Range("M88").Value = 223
$ws.Range("H91").Value = 183
$ws.Range("I91").Value = 183
$ws.Range("K91").Value = 183
$ws.Range("M91").Value = 1221
$ws.Range("H132").Value = 3846.8572
$ws.Range("I132").Value = 3972.6206
$ws.Range("J132").Value = 3239
$ws.Range("K132").Value = 11917.8618
$ws.Range("L132").Value = 9717
$ws.Range("M132").Value = -9387.861800000001
$ws.Range("N132").Value = -14777
